# edit.ps1 - applies "aggiunta di fase progettazione nel documento"
#
# Summary of changes:
#  1. Insert a new "Fase di progettazione" heading + two body paragraphs
#     at the very start of the document (before "Installazione").
#  2. Add <w:lastRenderedPageBreak/> before the run "Miriairim".
#  3. Remove <w:lastRenderedPageBreak/> before the run "Utenti del Sistema".
#  4. Insert two empty paragraphs right after the "Admin ha il compito..."
#     paragraph (before "Elenco delle funzionalità").
#  5. Remove <w:lastRenderedPageBreak/> before the run "rimuovere un gioco".
#  6. Add <w:lastRenderedPageBreak/> before the "Descrizione" run that
#     follows "Gestione Profilo Utente (Luca)".
#  7. Remove <w:lastRenderedPageBreak/> before the "Descrizione" run that
#     follows "Registrarsi alla piattaforma (Luca)".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert the three new paragraphs at the very beginning of the body.
# ---------------------------------------------------------------------
$introXml = '<w:p><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="4" w:color="EAECEF"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="360" w:after="240" w:line="240" w:lineRule="auto"/><w:outlineLvl w:val="0"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:kern w:val="36"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:kern w:val="36"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>Fase di progettazione</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="240" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Durante la fase di progettazione, abbiamo contribuito tutti allo stesso modo nella generazione dell’OOA e dell’OOD, nonché dello schema ER. Durante la progettazione, ci siamo resi conto che il sistema poteva essere diviso in 3 parti distinte non comunicanti.</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="240" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Pertanto, una volta finito, ognuno ha implementato una parte, apportando delle modifiche alla progettazione quando, durante l’implementazione, se ne è ritenuto necessario.</w:t></w:r></w:p>'

$startRng = $d.Range(0, 0)
$startRng.InsertXML($introXml)

# ---------------------------------------------------------------------
# 2. Add <w:lastRenderedPageBreak/> before "Miriairim".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Miriairim", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1).Range
$xml = '<w:p w:rsidR="00B614E3" w:rsidRDefault="00B614E3" w:rsidP="00B614E3"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Miriairim</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$para.InsertXML($xml)

# ---------------------------------------------------------------------
# 3. Remove <w:lastRenderedPageBreak/> before "Utenti del Sistema".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Utenti del Sistema", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1).Range
$xml = '<w:p w:rsidR="00C57DD1" w:rsidRPr="00C57DD1" w:rsidRDefault="00C57DD1" w:rsidP="00C57DD1"><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="4" w:color="EAECEF"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="360" w:after="240" w:line="240" w:lineRule="auto"/><w:outlineLvl w:val="0"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:kern w:val="36"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr><w:r w:rsidRPr="00C57DD1"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:kern w:val="36"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>Utenti del Sistema</w:t></w:r></w:p>'
$para.InsertXML($xml)

# ---------------------------------------------------------------------
# 4. Insert two empty paragraphs after "Admin ha il compito di gestire...".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("ha il compito di gestire il portale", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1).Range
$xml = '<w:p w:rsidR="00C57DD1" w:rsidRPr="00C57DD1" w:rsidRDefault="00C57DD1" w:rsidP="00C57DD1"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="60" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00C57DD1"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Admin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00C57DD1"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>ha il compito di gestire il portale. Tra le sue funzioni vi sono quelle di gestione degli utenti, gestione.</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="60" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="60" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
$para.InsertXML($xml)

# ---------------------------------------------------------------------
# 5. Remove <w:lastRenderedPageBreak/> before "rimuovere un gioco".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("rimuovere un gioco", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1).Range
$xml = '<w:p w:rsidR="0089427A" w:rsidRDefault="0089427A" w:rsidP="0089427A"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="60" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/></w:rPr><w:t>rimuovere un gioco</w:t></w:r></w:p>'
$para.InsertXML($xml)

# ---------------------------------------------------------------------
# 6. Add <w:lastRenderedPageBreak/> before the "Descrizione" run that
#    follows "Gestione Profilo Utente (Luca)".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Gestione Profilo Utente (Luca)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.Find.Execute("Descrizione", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1).Range
$xml = '<w:p w:rsidR="0089427A" w:rsidRDefault="0089427A" w:rsidP="0089427A"><w:pPr><w:pStyle w:val="Titolo3"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="360" w:after="240"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:lastRenderedPageBreak/><w:t>Descrizione</w:t></w:r></w:p>'
$para.InsertXML($xml)

# ---------------------------------------------------------------------
# 7. Remove <w:lastRenderedPageBreak/> before the "Descrizione" run that
#    follows "Registrarsi alla piattaforma (Luca)".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Registrarsi alla piattaforma (Luca)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.Find.Execute("Descrizione", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1).Range
$xml = '<w:p w:rsidR="0089427A" w:rsidRDefault="0089427A" w:rsidP="0089427A"><w:pPr><w:pStyle w:val="Titolo3"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="360" w:after="240"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>Descrizione</w:t></w:r></w:p>'
$para.InsertXML($xml)

Write-Host "All edits applied."
